# Saldo.xlsx update:
#  - Insert two new account rows (EDMUNDO / ANA-50k) right before the
#    "004581652 / CINCO" row (currently row 3).
#  - Remove the older duplicate "005198093 / ANA / 824.01" row further
#    down the sheet (it is superseded by the new ANA row added above).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert the two new rows above row 3 ("004581652 | CINCO | 14455.12") ---
$ws.Rows.Item(3).Resize(2).Insert()

# New row 3: account 003301389 / EDMUNDO / 50000
$ws.Cells.Item(3, 1).Value = "'003301389"
$ws.Cells.Item(3, 2).Value = "EDMUNDO"
$ws.Cells.Item(3, 3).Value = 50000

# New row 4: account 005198093 / ANA / 20824.01
$ws.Cells.Item(4, 1).Value = "'005198093"
$ws.Cells.Item(4, 2).Value = "ANA"
$ws.Cells.Item(4, 3).Value = 20824.01

# --- Remove the stale duplicate ANA row (824.01) that used to sit right
#     after "004482102 | NATALIA | 856.47". After the insert above it has
#     shifted down two rows, to row 16. ---
$ws.Rows.Item(16).Delete()
